$d = $word.ActiveDocument

# 1. Update the main heading text.
$d.Content.Find.Execute("Implante de Ressincronizador Convencional (CRT-P Serena™)", $true, $false, $false, $false, $false, $true, 1, $false, "CRT-P Serena", 2)

# 2. Remove the old subtitle paragraph entirely (including its paragraph mark).
$sub = $d.Paragraphs(2)
if ($sub.Range.Text.StartsWith("Terapia de ressincroniza")) {
    $sub.Range.Delete()
}

# 3. Update the material list entries: add a leading bullet glyph and trim/adjust wording.
$d.Content.Find.Execute("Gerador – Serena™ CRT-P", $true, $false, $false, $false, $false, $true, 1, $false, "• Gerador Serena", 2)
$d.Content.Find.Execute("Eletrodo Atrial – 5076-52", $true, $false, $false, $false, $false, $true, 1, $false, "• Eletrodo 5076-52", 2)
$d.Content.Find.Execute("Eletrodo Ventricular Direito – 5076-58", $true, $false, $false, $false, $false, $true, 1, $false, "• Eletrodo 5076-58", 2)
$d.Content.Find.Execute("Eletrodo Ventricular Esquerdo – 4298/4299", $true, $false, $false, $false, $false, $true, 1, $false, "• Eletrodo VE 4298/4299", 2)
$d.Content.Find.Execute("Bainha – 6250VIC", $true, $false, $false, $false, $false, $true, 1, $false, "• Bainha 6250VIC", 2)
$d.Content.Find.Execute("Ferramenta de Corte", $true, $false, $false, $false, $false, $true, 1, $false, "• Ferramenta de corte", 2)
$d.Content.Find.Execute("Guia 0.014", $true, $false, $false, $false, $false, $true, 1, $false, "• Guia 014", 2)
$d.Content.Find.Execute("Subseletora", $true, $false, $false, $false, $false, $true, 1, $false, "• Subseletora", 2)
$d.Content.Find.Execute("Introdutor – 3", $true, $false, $false, $false, $false, $true, 1, $false, "• Introdutor – 3", 2)
